$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared-string / label text ---
$ws.Range("A31").Value = "item26 [tonne*km]"
$ws.Range("A49").Value = "E_item [kWh]"

# --- Remove the duplicate bold/border style so A2:A26 & A31:A32 (previously
#     style index 2, a byte-for-byte duplicate of style index 1) collapse
#     onto the single remaining style. Re-apply the same visual style
#     (copy from a cell that already carries it, e.g. A1) to normalize.
$srcStyle = $ws.Range("A1")
$mergedRanges = @("A2:A3","A4:A5","A6:A8","A9:A13","A14:A15","A16:A19","A20:A21","A22:A24","A25:A26","A31:A32")
foreach ($rng in $mergedRanges) {
    $ws.Range($rng).Borders.LineStyle = $srcStyle.Borders.LineStyle
    $ws.Range($rng).Font.Bold = $srcStyle.Font.Bold
    $ws.Range($rng).HorizontalAlignment = $srcStyle.HorizontalAlignment
    $ws.Range($rng).VerticalAlignment = $srcStyle.VerticalAlignment
}

# --- Updated LCA values (new biosteam results) ---
$ws.Range("B37").Value = 4911701.166447819
$ws.Range("C37").Value = 137527632.6605389
$ws.Range("D37").Value = 1.206121793078792

$ws.Range("B38").Value = 2982.083297928013
$ws.Range("C38").Value = 790252.0739509233
$ws.Range("D38").Value = 0.006930536285537385

$ws.Range("D39").Value = -0.026044555768803
$ws.Range("C39").Value = -2969721.730535773

$ws.Range("C40").Value = -17438111.03414943
$ws.Range("D40").Value = -0.1529327986058634

$ws.Range("B41").Value = 345811.8509762709
$ws.Range("C41").Value = -1694478.069783727
$ws.Range("D41").Value = -0.01486062755769848

$ws.Range("B42").Value = 301871.7843221913
$ws.Range("D42").Value = -0.003971137989613406

$ws.Range("B43").Value = 127963.0998558381
$ws.Range("C43").Value = -691000.7392215255
$ws.Range("D43").Value = -0.006060098865119015

$ws.Range("B44").Value = 213694.1023588947
$ws.Range("D44").Value = -0.009183110577231602

$ws.Range("C45").Value = 114024664.3827575
